$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45175 to 45177
# (2023-09-06 -> 2023-09-08) for every data row (rows 2 through 422).
$ws.Range("C2:C422").Value = 45177
